# "每日学习.xlsx" — add the 2023-07-26 ("27") study-log entry as a new
# row 11 on Sheet1, and move the sheet's active-cell selection to E12
# (the cell just below the new row), matching the author's upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entry: day 27, two time ranges + two content notes.
$ws.Range("A11").Value = 27
$ws.Range("B11").Value = "10：24—11：27"
$ws.Range("C11").Value = "操作符"
$ws.Range("D11").Value = "3：15-7：20"
$ws.Range("E11").Value = "常见的一些关键字 常量和宏 指针 结构体 "

# Move the selection down to where the user's cursor ended up next.
$ws.Range("E12").Select()
